# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells that get refreshed each time
# the handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 66da0445-...md
$wsOverview.Range("G2").Value = "2016-08-26 11:02:45"

# zh-cn sheet, row 2 (66da0445-...zh-cn.xlf):
#   Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-08-26 11:02:40"
$wsZhCn.Range("K2").Value = "2016-08-26 11:03:08"

# de-de sheet, row 2 (66da0445-...de-de.xlf):
#   Correspond Handoff Datetime (H2) shares the same value/text as Overview!G2
#   Correspond Handback DateTime (K2)
$wsDeDe.Range("H2").Value = "2016-08-26 11:02:45"
$wsDeDe.Range("K2").Value = "2016-08-26 11:03:15"
